# Roll the GSC export window forward by one day:
#   - drop the oldest date row (2025-10-10) from the "Chart" sheet
#   - every remaining row shifts up by one
#   - the sheet now ends one row earlier (A1:C91 -> A1:C90)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Rows.Item(2).Delete()
